$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 29.12242425684365
$ws.Range("C2").Value = 23.02353161545612
$ws.Range("D2").Value = 5.457562313420514
$ws.Range("E2").Value = 29.30942974478529
$ws.Range("F2").Value = 42.83252900737337
$ws.Range("G2").Value = 2.069037592771776
$ws.Range("H2").Value = 3.219152413621614
$ws.Range("I2").Value = 3.413434991520969
$ws.Range("P2").Value = 13.37876734696289
$ws.Range("B3").Value = 27.1374243989585
$ws.Range("C3").Value = 21.40733941091886
$ws.Range("D3").Value = 5.337267986424786
$ws.Range("E3").Value = 27.25656006138376
$ws.Range("F3").Value = 40.45473511116975
$ws.Range("G3").Value = 2.078198967281577
$ws.Range("H3").Value = 2.859703050369023
$ws.Range("I3").Value = 3.074111969351026
$ws.Range("P3").Value = 13.44169929310584
$ws.Range("B4").Value = 25.87105930567509
$ws.Range("C4").Value = 20.36552448339359
$ws.Range("D4").Value = 5.260392023805906
$ws.Range("E4").Value = 25.9310335883985
$ws.Range("F4").Value = 38.9386664220199
$ws.Range("G4").Value = 2.083948925629955
$ws.Range("H4").Value = 2.633168838164859
$ws.Range("I4").Value = 2.861811536839314
$ws.Range("P4").Value = 13.48005957020183
$ws.Range("B5").Value = 25.32578769390431
$ws.Range("C5").Value = 19.93426239784606
$ws.Range("D5").Value = 5.224182377235195
$ws.Range("E5").Value = 25.37362898793707
$ws.Range("F5").Value = 38.28229444321997
$ws.Range("G5").Value = 2.086343736601319
$ws.Range("H5").Value = 2.538891211510318
$ws.Range("I5").Value = 2.774232242323455
$ws.Range("P5").Value = 13.49076358101514
$ws.Range("B6").Value = 25.21926499381419
$ws.Range("C6").Value = 19.86909335459384
$ws.Range("D6").Value = 5.213083032900629
$ws.Range("E6").Value = 25.27925579919151
$ws.Range("F6").Value = 38.14214131045867
$ws.Range("G6").Value = 2.086765773327855
$ws.Range("H6").Value = 2.522620236960007
$ws.Range("I6").Value = 2.759813657969449
$ws.Range("P6").Value = 13.48663518098441
$ws.Range("B7").Value = 25.82386594742595
$ws.Range("C7").Value = 20.37891803677698
$ws.Range("D7").Value = 5.24621841167165
$ws.Range("E7").Value = 25.92137325416789
$ws.Range("F7").Value = 38.84853422600711
$ws.Range("G7").Value = 2.084041586545941
$ws.Range("H7").Value = 2.630487630487425
$ws.Range("I7").Value = 2.860669597260886
$ws.Range("P7").Value = 13.46411645950251
$ws.Range("B8").Value = 28.39101441007527
$ws.Range("C8").Value = 22.49958675012746
$ws.Range("D8").Value = 5.398999121240402
$ws.Range("E8").Value = 28.61202614888293
$ws.Range("F8").Value = 41.92355994330617
$ws.Range("G8").Value = 2.072250082095517
$ws.Range("H8").Value = 3.094540212180184
$ws.Range("I8").Value = 3.296299062931547
$ws.Range("P8").Value = 13.3793457196131
$ws.Range("B9").Value = 33.0895634734422
$ws.Range("C9").Value = 26.24400022884295
$ws.Range("D9").Value = 5.698493296056921
$ws.Range("E9").Value = 33.386799994948
$ws.Range("F9").Value = 47.63182860177298
$ws.Range("G9").Value = 2.049914329690132
$ws.Range("H9").Value = 3.964171270026597
$ws.Range("I9").Value = 4.12457982399691
$ws.Range("P9").Value = 13.24067711716932
$ws.Range("B10").Value = 36.0434723756846
$ws.Range("C10").Value = 28.63441750164197
$ws.Range("D10").Value = 5.830995088903073
$ws.Range("E10").Value = 35.71059044858454
$ws.Range("F10").Value = 51.13791178948935
$ws.Range("G10").Value = 2.034521241778429
$ws.Range("H10").Value = 4.534438610626464
$ws.Range("I10").Value = 4.698664340282644
$ws.Range("P10").Value = 13.06500675699586
$ws.Range("B11").Value = 36.02269999258484
$ws.Range("C11").Value = 28.45708531374927
$ws.Range("D11").Value = 5.314150617902084
$ws.Range("E11").Value = 29.16197763703897
$ws.Range("F11").Value = 49.38345475826942
$ws.Range("G11").Value = 2.03254108020656
$ws.Range("H11").Value = 4.85817995592064
$ws.Range("I11").Value = 4.774015261544365
$ws.Range("P11").Value = 12.41367954742661
$ws.Range("B12").Value = 35.45132614581257
$ws.Range("C12").Value = 27.80348711356223
$ws.Range("D12").Value = 4.902634528648329
$ws.Range("E12").Value = 23.22781200176598
$ws.Range("F12").Value = 47.2578856750118
$ws.Range("G12").Value = 2.033596121678061
$ws.Range("H12").Value = 5.582613299285839
$ws.Range("I12").Value = 4.732397001373103
$ws.Range("P12").Value = 11.97410097664093
$ws.Range("B13").Value = 34.36056989042019
$ws.Range("C13").Value = 26.74867764636883
$ws.Range("D13").Value = 4.535521992404589
$ws.Range("E13").Value = 17.30423646035947
$ws.Range("F13").Value = 44.55714983913662
$ws.Range("G13").Value = 2.037103107177901
$ws.Range("H13").Value = 6.5172900259235
$ws.Range("I13").Value = 4.598145302969217
$ws.Range("P13").Value = 11.65383825406395
$ws.Range("B14").Value = 33.34101321820851
$ws.Range("C14").Value = 25.82294158201565
$ws.Range("D14").Value = 4.309114174977028
$ws.Range("E14").Value = 13.23754183903135
$ws.Range("F14").Value = 42.3561915673402
$ws.Range("G14").Value = 2.040608279471432
$ws.Range("H14").Value = 7.264252431868699
$ws.Range("I14").Value = 4.464898670560597
$ws.Range("P14").Value = 11.49213779797367
$ws.Range("B15").Value = 32.9648862501232
$ws.Range("C15").Value = 25.51230040168827
$ws.Range("D15").Value = 4.257088291712681
$ws.Range("E15").Value = 12.24706028960986
$ws.Range("F15").Value = 41.66581721299364
$ws.Range("G15").Value = 2.042053549585688
$ws.Range("H15").Value = 7.433725654074647
$ws.Range("I15").Value = 4.411213507189735
$ws.Range("P15").Value = 11.46914446365359
$ws.Range("B16").Value = 31.87608673764142
$ws.Range("C16").Value = 24.67704261173039
$ws.Range("D16").Value = 4.270210178006453
$ws.Range("E16").Value = 11.94105676579287
$ws.Range("F16").Value = 40.52409631443208
$ws.Range("G16").Value = 2.047802731150009
$ws.Range("H16").Value = 7.12150906140563
$ws.Range("I16").Value = 4.199051954891138
$ws.Range("P16").Value = 11.6109875496389
$ws.Range("B17").Value = 31.59081828550146
$ws.Range("C17").Value = 24.54023630133161
$ws.Range("D17").Value = 4.392134655494539
$ws.Range("E17").Value = 13.92327043899569
$ws.Range("F17").Value = 40.8602284027189
$ws.Range("G17").Value = 2.050415076146294
$ws.Range("H17").Value = 6.426382675262424
$ws.Range("I17").Value = 4.104955505081555
$ws.Range("P17").Value = 11.7975146447877
$ws.Range("B18").Value = 32.01266221312296
$ws.Range("C18").Value = 24.99574270966277
$ws.Range("D18").Value = 4.659682672957626
$ws.Range("E18").Value = 18.43728141267208
$ws.Range("F18").Value = 42.57297928984077
$ws.Range("G18").Value = 2.050369575357616
$ws.Range("H18").Value = 5.417158518650544
$ws.Range("I18").Value = 4.10792768057442
$ws.Range("P18").Value = 12.08228754755286
$ws.Range("B19").Value = 32.91428265727141
$ws.Range("C19").Value = 25.92709831785335
$ws.Range("D19").Value = 5.046921090650943
$ws.Range("E19").Value = 24.73815109591906
$ws.Range("F19").Value = 45.11864271034842
$ws.Range("G19").Value = 2.04795465780415
$ws.Range("H19").Value = 4.505350713937855
$ws.Range("I19").Value = 4.201341048233777
$ws.Range("P19").Value = 12.44996205798263
$ws.Range("B20").Value = 35.18532413088438
$ws.Range("C20").Value = 28.05782744234055
$ws.Range("D20").Value = 5.75501748214636
$ws.Range("E20").Value = 35.0635044301351
$ws.Range("F20").Value = 50.01271585242392
$ws.Range("G20").Value = 2.038744722034453
$ws.Range("H20").Value = 4.378823149754055
$ws.Range("I20").Value = 4.545497348120654
$ws.Range("P20").Value = 13.05914361158725
$ws.Range("B21").Value = 37.60410681397683
$ws.Range("C21").Value = 30.04120099125489
$ws.Range("D21").Value = 5.965210850010251
$ws.Range("E21").Value = 38.14978486046395
$ws.Range("F21").Value = 53.24885792036237
$ws.Range("G21").Value = 2.026048842488905
$ws.Range("H21").Value = 4.886079521493445
$ws.Range("I21").Value = 5.01815482679849
$ws.Range("P21").Value = 13.03976532144044
$ws.Range("B22").Value = 39.08938018609194
$ws.Range("C22").Value = 31.20984196780294
$ws.Range("D22").Value = 6.081779563242853
$ws.Range("E22").Value = 39.6612641198891
$ws.Range("F22").Value = 55.20015203280625
$ws.Range("G22").Value = 2.01799549018383
$ws.Range("H22").Value = 5.192950727496313
$ws.Range("I22").Value = 5.317053975957242
$ws.Range("P22").Value = 13.01478655940469
$ws.Range("B23").Value = 38.3363028547522
$ws.Range("C23").Value = 30.57376622537857
$ws.Range("D23").Value = 6.034294527975277
$ws.Range("E23").Value = 38.86161477420278
$ws.Range("F23").Value = 54.2361464001029
$ws.Range("G23").Value = 2.022226197195321
$ws.Range("H23").Value = 5.030911046414971
$ws.Range("I23").Value = 5.158508503327936
$ws.Range("P23").Value = 13.04675674622192
$ws.Range("B24").Value = 35.31861004282749
$ws.Range("C24").Value = 28.11630933623611
$ws.Range("D24").Value = 5.827425069724348
$ws.Range("E24").Value = 35.73341028950765
$ws.Range("F24").Value = 50.37637003397209
$ws.Range("G24").Value = 2.038380818390848
$ws.Range("H24").Value = 4.411048834133488
$ws.Range("I24").Value = 4.555896700817713
$ws.Range("P24").Value = 13.14098249754714
$ws.Range("B25").Value = 31.81816931407905
$ws.Range("C25").Value = 25.29905695135953
$ws.Range("D25").Value = 5.596600415070315
$ws.Range("E25").Value = 32.14356109270349
$ws.Range("F25").Value = 46.01060943009153
$ws.Range("G25").Value = 2.055980476337377
$ws.Range("H25").Value = 3.730515494968214
$ws.Range("I25").Value = 3.901720770525965
$ws.Range("P25").Value = 13.2493814597989
